$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 13) mirroring the existing rows' structure.
$row = 13

# Copy the formatting (including the date number format) from the row above
# so the new cell reuses the existing style index instead of creating a new one.
$ws.Range("A12").Copy()
$ws.Range("A$row").PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = 42620.891203703701

$ws.Cells.Item($row, 2).Value = 2
$ws.Cells.Item($row, 3).Value = 55
$ws.Cells.Item($row, 4).Value = 41
$ws.Cells.Item($row, 5).Value = 55
$ws.Cells.Item($row, 6).Value = 50
$ws.Cells.Item($row, 7).Value = 31857
$ws.Cells.Item($row, 8).Value = 28453
$ws.Cells.Item($row, 9).Value = 3144
$ws.Cells.Item($row, 10).Value = 461
$ws.Cells.Item($row, 11).Value = 347
$ws.Cells.Item($row, 12).Value = 2
$ws.Cells.Item($row, 13).Value = 2
$ws.Cells.Item($row, 14).Value = "Bag"
